$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.904.01'
$ws.Range('E2').Value = '  +0.64%  '
$ws.Range('D3').Value = '2.812.19'
$ws.Range('E3').Value = '  +1.53%  '
$ws.Range('D5').Value = "'357.01"
$ws.Range('E5').Value = '  +0.26%  '
$ws.Range('D6').Value = "'112.57"
$ws.Range('E6').Value = '  +3.45%  '
$ws.Range('D7').Value = "'0.557"
$ws.Range('E7').Value = '  +0.59%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = "'0.632"
$ws.Range('D10').Value = "'40.47"
$ws.Range('E10').Value = '  +2.63%  '
$ws.Range('E11').Value = '  -0.41%  '
$ws.Range('E12').Value = '  +0.17%  '
$ws.Range('D13').Value = "'20.01"
$ws.Range('E13').Value = '  +3.30%  '
$ws.Range('D14').Value = "'7.80"
$ws.Range('E14').Value = '  +3.23%  '
$ws.Range('D15').Value = '3.259.24'
$ws.Range('E15').Value = '  +1.54%  '
$ws.Range('D16').Value = '2.818.78'
$ws.Range('E16').Value = '  +2.18%  '
$ws.Range('D17').Value = "'0.944"
$ws.Range('E17').Value = '  +1.53%  '
$ws.Range('D18').Value = '51.917.06'
$ws.Range('E18').Value = '  +0.73%  '
$ws.Range('D19').Value = "'7.65"
$ws.Range('E19').Value = '  +3.30%  '
$ws.Range('E20').Value = '  +3.32%  '
$ws.Range('D21').Value = "'13.64"
$ws.Range('E21').Value = '  +4.48%  '
$ws.Range('D22').Value = '0.0₃0982'
$ws.Range('E22').Value = '  +1.85%  '
$ws.Range('D23').Value = "'70.48"
$ws.Range('E23').Value = '  +0.87%  '
$ws.Range('D24').Value = "'268.82"
$ws.Range('E24').Value = '  +0.64%  '
$ws.Range('D25').Value = "'2.78"
$ws.Range('E25').Value = '  +1.75%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = "'1.00"
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = "'26.19"
$ws.Range('E27').Value = '  -0.43%  '
$ws.Range('E28').Value = '  -0.12%  '
$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D29').Value = "'38.42"
$ws.Range('E29').Value = '  +11.65%  '
$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').Value = "'10.43"
$ws.Range('E30').Value = '  +2.80%  '
$ws.Range('D31').Value = "'2.24"
$ws.Range('E31').Value = '  +1.48%  '
$ws.Range('D32').Value = "'6.19"
$ws.Range('E32').Value = '  +0.25%  '
$ws.Range('D33').Value = "'52.48"
$ws.Range('E33').Value = '  +1.69%  '
$ws.Range('D34').Value = "'5.65"
$ws.Range('E34').Value = '  +9.99%  '
$ws.Range('E35').Value = '  -0.72%  '
$ws.Range('D36').Value = "'0.0879"
$ws.Range('E36').Value = '  +5.29%  '
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('D38').Value = "'18.89"
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('D39').Value = "'2.01"
$ws.Range('E39').Value = '  +3.50%  '
$ws.Range('D40').Value = "'3.14"
$ws.Range('E40').Value = '  +1.05%  '
$ws.Range('E41').Value = '  +1.46%  '
$ws.Range('D42').Value = "'2.52"
$ws.Range('E42').Value = '  -0.45%  '
$ws.Range('D43').Value = "'120.92"
$ws.Range('E43').Value = '  +1.34%  '
$ws.Range('D44').Value = "'21.98"
$ws.Range('E44').Value = '  +1.98%  '
$ws.Range('E45').Value = '  -0.96%  '
$ws.Range('D46').Value = "'3.40"
$ws.Range('E46').Value = '  +4.88%  '
$ws.Range('D47').Value = '2.110.85'
$ws.Range('E47').Value = '  +1.29%  '
$ws.Range('D48').Value = "'2.40"
$ws.Range('E48').Value = '  +5.19%  '
$ws.Range('D49').Value = "'0.944"
$ws.Range('E49').Value = '  +1.34%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').Value = "'5.47"
$ws.Range('E50').Value = '  -1.18%  '
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').Value = "'1.36"
$ws.Range('E51').Value = '  +7.83%  '

# Reset style on text-forced numeric-looking cells so no extra style index lingers
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
